$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple Porcentaje (column I) value updates
$ws.Range("I15").Value = 0.25
$ws.Range("I16").Value = 0.2
$ws.Range("I61").Value = 0.4
$ws.Range("I62").Value = 0.6
$ws.Range("I125").Value = 0.15
$ws.Range("I126").Value = 0.15
$ws.Range("I127").Value = 0.15
$ws.Range("I128").Value = 0.15
$ws.Range("I129").Value = 0.25
$ws.Range("I130").Value = 0.35

# Row 131: sopa / fideo / cereal / crema / #FFFACD / 0.35
$ws.Range("E131").Value = "fideo"
$ws.Range("F131").Value = "cereal"
$ws.Range("G131").Value = "crema"
$ws.Range("H131").Value = "#FFFACD"
$ws.Range("H131").Interior.Color = 13499135
$ws.Range("I131").Value = 0.35

# Row 132: Y / sopa / queso / lacteo / crema / #FFFACD / 0.05
$ws.Range("C132").Value = "Y"
$ws.Range("E132").Value = "queso"
$ws.Range("F132").Value = "lacteo"
$ws.Range("G132").Value = "crema"
$ws.Range("H132").Value = "#FFFACD"
$ws.Range("H132").Interior.Color = 13499135
$ws.Range("I132").Value = 0.05

# Delete rows 133 and 134 (shift remaining rows up)
$ws.Rows.Item(134).Delete()
$ws.Rows.Item(133).Delete()
